$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "MEC-3A-Trat. Térmicos"
$ws.Range("B12").Value = "MEC-3A-Trat. Térmicos"
$ws.Range("D12").Value = "-"
$ws.Range("B14").Value = "MEC-3A-Trat. Térmicos"
$ws.Range("D14").Value = "-"
$ws.Range("B15").Value = "MEC-3A-Trat. Térmicos"
$ws.Range("D15").Value = "-"
$ws.Range("D16").Value = "-"
